$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Locate the paragraphs we need to touch by their current text instead of a
# hard-coded index, so the script is resilient to how the collection is
# walked.
$pGoBack = $null
$pShip = $null
$pEngine = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    # Range.Text carries the trailing paragraph mark (CR); strip it before
    # comparing against the plain paragraph text.
    $t = $cand.Range.Text.TrimEnd("`r", "`n", "`a")
    if ($t -eq "Ship") {
        $pShip = $cand
    } elseif ($t.StartsWith("Engine") -and $t.Contains("flame")) {
        $pEngine = $cand
    } elseif ($t -eq "Interface Sketch") {
        # The "_GoBack" bookmark paragraph sits immediately after the
        # "Interface Sketch" heading and carries no visible text.
        $pGoBack = $d.Paragraphs.Item($i + 1)
    }
}

# ---------------------------------------------------------------------------
# 1) The paragraph right after "Interface Sketch" used to hold only the
#    "_GoBack" bookmark. The bookmark is removed, leaving an empty paragraph.
# ---------------------------------------------------------------------------
$rGoBack = $pGoBack.Range
$rGoBack.MoveEnd(1, -1) | Out-Null
$rGoBack.InsertXML("<w:p $wns></w:p>") | Out-Null

# ---------------------------------------------------------------------------
# 2) The "Ship" paragraph gains extra runs describing where the asset came
#    from: " - From Free SciFi Fighter Pubvlisher CGPitbull"
# ---------------------------------------------------------------------------
$rShip = $pShip.Range
$rShip.MoveEnd(1, -1) | Out-Null
$shipXml = "<w:p $wns>" +
  "<w:r><w:t xml:space=`"preserve`"> – From Free </w:t></w:r>" +
  "<w:proofErr w:type=`"spellStart`"/>" +
  "<w:r><w:t>SciFi</w:t></w:r>" +
  "<w:proofErr w:type=`"spellEnd`"/>" +
  "<w:r><w:t xml:space=`"preserve`"> Fighter </w:t></w:r>" +
  "<w:proofErr w:type=`"spellStart`"/>" +
  "<w:r><w:t>Pubvlisher</w:t></w:r>" +
  "<w:proofErr w:type=`"spellEnd`"/>" +
  "<w:r><w:t xml:space=`"preserve`"> </w:t></w:r>" +
  "<w:proofErr w:type=`"spellStart`"/>" +
  "<w:r><w:t>CGPitbull</w:t></w:r>" +
  "<w:proofErr w:type=`"spellEnd`"/>" +
  "</w:p>"
$rShip.InsertXML($shipXml) | Out-Null

# ---------------------------------------------------------------------------
# 3) The old "Engine - flame 1 , flame 2" paragraph is replaced wholesale by
#    two new paragraphs describing the engine particle system assets.
# ---------------------------------------------------------------------------
$rEngine = $pEngine.Range
$engineXml = "<w:p $wns>" +
  "<w:r><w:t>Engine –</w:t></w:r>" +
  "<w:r><w:t xml:space=`"preserve`"> From </w:t></w:r>" +
  "<w:proofErr w:type=`"spellStart`"/>" +
  "<w:r><w:t>SimpleParticlePack</w:t></w:r>" +
  "<w:proofErr w:type=`"spellEnd`"/>" +
  "<w:r><w:t xml:space=`"preserve`"> Publisher Unity Technologies</w:t></w:r>" +
  "</w:p>" +
  "<w:p $wns>" +
  "<w:r><w:t xml:space=`"preserve`"> </w:t></w:r>" +
  "<w:proofErr w:type=`"spellStart`"/>" +
  "<w:r><w:t>SimpleFlame</w:t></w:r>" +
  "<w:proofErr w:type=`"spellEnd`"/>" +
  "<w:r><w:t xml:space=`"preserve`">(Blue), </w:t></w:r>" +
  "<w:proofErr w:type=`"spellStart`"/>" +
  "<w:r><w:t>SimpleFlame</w:t></w:r>" +
  "<w:proofErr w:type=`"spellEnd`"/>" +
  "<w:r><w:t xml:space=`"preserve`">(Green), </w:t></w:r>" +
  "<w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/>" +
  "<w:bookmarkEnd w:id=`"0`"/>" +
  "<w:r><w:t>Torch(Green) x2</w:t></w:r>" +
  "</w:p>"
$rEngine.InsertXML($engineXml) | Out-Null
